# test cases 3 and 4
# Adds two new worksheets (NumOfFavorites, OneFavorite) representing the
# "apartment favorites" test fixtures, updates the DetailedSearch sheet's
# listing title, and makes OneFavorite the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # BasicSearch
$ws2 = $wb.Worksheets.Item(2)   # DetailedSearch

# ---------------------------------------------------------------------
# 1. Create "NumOfFavorites" by duplicating DetailedSearch (keeps the
#    same column-A width / fonts / number formats as a starting point).
# ---------------------------------------------------------------------
$ws2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "NumOfFavorites"

# ---------------------------------------------------------------------
# 2. Create "OneFavorite" by duplicating BasicSearch (no custom column
#    widths, matching the target layout) and becomes the active sheet.
# ---------------------------------------------------------------------
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "OneFavorite"

# ---------------------------------------------------------------------
# 3. DetailedSearch: rename the listing and widen column B.
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "Halletts Point"
$ws2.Range("B2").Font.Name = "Menlo"
$ws2.Range("B2").Font.Size = 11
$ws2.Columns.Item(2).ColumnWidth = 22.666666666666668

# ---------------------------------------------------------------------
# 4. NumOfFavorites: Queens listing with its favorites count.
# ---------------------------------------------------------------------
$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A2").Value = "Queens, NY"

$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "2"
$ws3.Range("B2").Font.Name = "Menlo"
$ws3.Range("B2").Font.Size = 11

$ws3.Columns.Item(1).ColumnWidth = 12.5
$ws3.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. OneFavorite: Queens listing with a single favorite.
# ---------------------------------------------------------------------
$ws4.Range("A2").Value = "Queens, NY"

$ws4.Range("B2").NumberFormat = "@"
$ws4.Range("B2").Value = "1"

$ws4.Range("E7").Select() | Out-Null

# Make OneFavorite the active tab (last sheet).
$ws4.Activate() | Out-Null
